# "Excel: added super simple expression support to filters for string
# matching."
#
# Populate the (previously empty) filter-expression column G, for the
# data rows (4-10), with the string "beer_kettle_01.jpg" so the filter
# can match against that file name. The new cells should look like the
# existing string column D (text-formatted, i.e. NumberFormat "@"),
# so copy that format over before writing the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textFormat = $ws.Range("D4").NumberFormat

for ($row = 4; $row -le 10; $row++) {
    $cell = $ws.Range("G$row")
    $cell.NumberFormat = $textFormat
    $cell.Value = "beer_kettle_01.jpg"
}
